$d = $word.ActiveDocument

$d.Content.Find.Execute("77-26=", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=", 2) | Out-Null
$d.Content.Find.Execute("93-49=", $true, $false, $false, $false, $false, $true, 1, $false, "53+15=", 2) | Out-Null
$d.Content.Find.Execute("84-57=", $true, $false, $false, $false, $false, $true, 1, $false, "51+46=", 2) | Out-Null
$d.Content.Find.Execute("69+4=", $true, $false, $false, $false, $false, $true, 1, $false, "95-1=", 2) | Out-Null
$d.Content.Find.Execute("30-11=", $true, $false, $false, $false, $false, $true, 1, $false, "37-21=", 2) | Out-Null
$d.Content.Find.Execute("92-27=", $true, $false, $false, $false, $false, $true, 1, $false, "48+2=", 2) | Out-Null
$d.Content.Find.Execute("43+48=", $true, $false, $false, $false, $false, $true, 1, $false, "45+14=", 2) | Out-Null
$d.Content.Find.Execute("43-16=", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=", 2) | Out-Null
$d.Content.Find.Execute("21+31=", $true, $false, $false, $false, $false, $true, 1, $false, "31+68=", 2) | Out-Null
$d.Content.Find.Execute("5+48=", $true, $false, $false, $false, $false, $true, 1, $false, "64-23=", 2) | Out-Null
$d.Content.Find.Execute("67+17=", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=", 2) | Out-Null
$d.Content.Find.Execute("97-49=", $true, $false, $false, $false, $false, $true, 1, $false, "91-81=", 2) | Out-Null
$d.Content.Find.Execute("9+54=", $true, $false, $false, $false, $false, $true, 1, $false, "73-46=", 2) | Out-Null
$d.Content.Find.Execute("10+10=", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=", 2) | Out-Null
$d.Content.Find.Execute("38-7=", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=", 2) | Out-Null
$d.Content.Find.Execute("1+38=", $true, $false, $false, $false, $false, $true, 1, $false, "97-74=", 2) | Out-Null
$d.Content.Find.Execute("98-42=", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=", 2) | Out-Null
$d.Content.Find.Execute("82-77=", $true, $false, $false, $false, $false, $true, 1, $false, "24+13=", 2) | Out-Null
$d.Content.Find.Execute("40+41=", $true, $false, $false, $false, $false, $true, 1, $false, "74-20=", 2) | Out-Null
$d.Content.Find.Execute("0+19=", $true, $false, $false, $false, $false, $true, 1, $false, "20+22=", 2) | Out-Null
$d.Content.Find.Execute("43+23=", $true, $false, $false, $false, $false, $true, 1, $false, "63+33=", 2) | Out-Null
$d.Content.Find.Execute("37+38=", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=", 2) | Out-Null
$d.Content.Find.Execute("23-4=", $true, $false, $false, $false, $false, $true, 1, $false, "38+55=", 2) | Out-Null
$d.Content.Find.Execute("61-9=", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=", 2) | Out-Null
$d.Content.Find.Execute("38+41=", $true, $false, $false, $false, $false, $true, 1, $false, "60+34=", 2) | Out-Null
$d.Content.Find.Execute("91-63=", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=", 2) | Out-Null
$d.Content.Find.Execute("8+27=", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=", 2) | Out-Null
$d.Content.Find.Execute("18-14=", $true, $false, $false, $false, $false, $true, 1, $false, "46+15=", 2) | Out-Null
$d.Content.Find.Execute("40-0=", $true, $false, $false, $false, $false, $true, 1, $false, "68-64=", 2) | Out-Null
$d.Content.Find.Execute("44-33=", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=", 2) | Out-Null
$d.Content.Find.Execute("15+53=", $true, $false, $false, $false, $false, $true, 1, $false, "92-82=", 2) | Out-Null
$d.Content.Find.Execute("92-44=", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=", 2) | Out-Null
$d.Content.Find.Execute("43-8=", $true, $false, $false, $false, $false, $true, 1, $false, "47+33=", 2) | Out-Null
$d.Content.Find.Execute("12+44=", $true, $false, $false, $false, $false, $true, 1, $false, "64-5=", 2) | Out-Null
$d.Content.Find.Execute("19-16=", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=", 2) | Out-Null
$d.Content.Find.Execute("44+53=", $true, $false, $false, $false, $false, $true, 1, $false, "6+74=", 2) | Out-Null
$d.Content.Find.Execute("88-76=", $true, $false, $false, $false, $false, $true, 1, $false, "73-31=", 2) | Out-Null
$d.Content.Find.Execute("41+13=", $true, $false, $false, $false, $false, $true, 1, $false, "72+16=", 2) | Out-Null
$d.Content.Find.Execute("52-24=", $true, $false, $false, $false, $false, $true, 1, $false, "19-13=", 2) | Out-Null
$d.Content.Find.Execute("1+94=", $true, $false, $false, $false, $false, $true, 1, $false, "57-40=", 2) | Out-Null
$d.Content.Find.Execute("52-31=", $true, $false, $false, $false, $false, $true, 1, $false, "45-38=", 2) | Out-Null
$d.Content.Find.Execute("8+60=", $true, $false, $false, $false, $false, $true, 1, $false, "42+33=", 2) | Out-Null
$d.Content.Find.Execute("36+40=", $true, $false, $false, $false, $false, $true, 1, $false, "66+28=", 2) | Out-Null
$d.Content.Find.Execute("92-16=", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=", 2) | Out-Null
$d.Content.Find.Execute("91+7=", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=", 2) | Out-Null
$d.Content.Find.Execute("68-21=", $true, $false, $false, $false, $false, $true, 1, $false, "72-5=", 2) | Out-Null
$d.Content.Find.Execute("70-12=", $true, $false, $false, $false, $false, $true, 1, $false, "58-14=", 2) | Out-Null
$d.Content.Find.Execute("73+21=", $true, $false, $false, $false, $false, $true, 1, $false, "33-17=", 2) | Out-Null
$d.Content.Find.Execute("93-92=", $true, $false, $false, $false, $false, $true, 1, $false, "66-15=", 2) | Out-Null
$d.Content.Find.Execute("45+49=", $true, $false, $false, $false, $false, $true, 1, $false, "45-27=", 2) | Out-Null
$d.Content.Find.Execute("39+31=", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=", 2) | Out-Null
$d.Content.Find.Execute("88-72=", $true, $false, $false, $false, $false, $true, 1, $false, "17+34=", 2) | Out-Null
$d.Content.Find.Execute("33+45=", $true, $false, $false, $false, $false, $true, 1, $false, "56+3=", 2) | Out-Null
$d.Content.Find.Execute("9+23=", $true, $false, $false, $false, $false, $true, 1, $false, "74-69=", 2) | Out-Null
$d.Content.Find.Execute("15+60=", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=", 2) | Out-Null
$d.Content.Find.Execute("20+8=", $true, $false, $false, $false, $false, $true, 1, $false, "35+26=", 2) | Out-Null
$d.Content.Find.Execute("83+13=", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=", 2) | Out-Null
$d.Content.Find.Execute("42+39=", $true, $false, $false, $false, $false, $true, 1, $false, "59+28=", 2) | Out-Null
$d.Content.Find.Execute("60-45=", $true, $false, $false, $false, $false, $true, 1, $false, "50-48=", 2) | Out-Null
$d.Content.Find.Execute("33+48=", $true, $false, $false, $false, $false, $true, 1, $false, "29+2=", 2) | Out-Null
$d.Content.Find.Execute("44-1=", $true, $false, $false, $false, $false, $true, 1, $false, "28-17=", 2) | Out-Null
$d.Content.Find.Execute("28+28=", $true, $false, $false, $false, $false, $true, 1, $false, "40+59=", 2) | Out-Null
$d.Content.Find.Execute("37+14=", $true, $false, $false, $false, $false, $true, 1, $false, "31+46=", 2) | Out-Null
$d.Content.Find.Execute("50+15=", $true, $false, $false, $false, $false, $true, 1, $false, "44-11=", 2) | Out-Null
$d.Content.Find.Execute("19+25=", $true, $false, $false, $false, $false, $true, 1, $false, "11+10=", 2) | Out-Null
$d.Content.Find.Execute("11+37=", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=", 2) | Out-Null
$d.Content.Find.Execute("75-44=", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=", 2) | Out-Null
$d.Content.Find.Execute("98-95=", $true, $false, $false, $false, $false, $true, 1, $false, "57+30=", 2) | Out-Null
$d.Content.Find.Execute("67-59=", $true, $false, $false, $false, $false, $true, 1, $false, "16+82=", 2) | Out-Null
$d.Content.Find.Execute("14+59=", $true, $false, $false, $false, $false, $true, 1, $false, "78+10=", 2) | Out-Null
$d.Content.Find.Execute("43-5=", $true, $false, $false, $false, $false, $true, 1, $false, "50-21=", 2) | Out-Null
$d.Content.Find.Execute("31+50=", $true, $false, $false, $false, $false, $true, 1, $false, "75-48=", 2) | Out-Null
$d.Content.Find.Execute("75-1=", $true, $false, $false, $false, $false, $true, 1, $false, "85-48=", 2) | Out-Null
$d.Content.Find.Execute("74+9=", $true, $false, $false, $false, $false, $true, 1, $false, "32+31=", 2) | Out-Null
$d.Content.Find.Execute("58-46=", $true, $false, $false, $false, $false, $true, 1, $false, "11+15=", 2) | Out-Null
$d.Content.Find.Execute("57+9=", $true, $false, $false, $false, $false, $true, 1, $false, "45-39=", 2) | Out-Null
$d.Content.Find.Execute("75-47=", $true, $false, $false, $false, $false, $true, 1, $false, "60-20=", 2) | Out-Null
$d.Content.Find.Execute("97-68=", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=", 2) | Out-Null
$d.Content.Find.Execute("6+90=", $true, $false, $false, $false, $false, $true, 1, $false, "17+36=", 2) | Out-Null
$d.Content.Find.Execute("63-59=", $true, $false, $false, $false, $false, $true, 1, $false, "5+63=", 2) | Out-Null
$d.Content.Find.Execute("92-31=", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=", 2) | Out-Null
$d.Content.Find.Execute("5+46=", $true, $false, $false, $false, $false, $true, 1, $false, "99-10=", 2) | Out-Null
$d.Content.Find.Execute("88-71=", $true, $false, $false, $false, $false, $true, 1, $false, "57-14=", 2) | Out-Null
$d.Content.Find.Execute("22+9=", $true, $false, $false, $false, $false, $true, 1, $false, "99-26=", 2) | Out-Null
$d.Content.Find.Execute("11+43=", $true, $false, $false, $false, $false, $true, 1, $false, "34+8=", 2) | Out-Null
$d.Content.Find.Execute("93-17=", $true, $false, $false, $false, $false, $true, 1, $false, "73-24=", 2) | Out-Null
$d.Content.Find.Execute("17-10=", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=", 2) | Out-Null
$d.Content.Find.Execute("39+36=", $true, $false, $false, $false, $false, $true, 1, $false, "88-47=", 2) | Out-Null
$d.Content.Find.Execute("96-60=", $true, $false, $false, $false, $false, $true, 1, $false, "63+14=", 2) | Out-Null
$d.Content.Find.Execute("40+7=", $true, $false, $false, $false, $false, $true, 1, $false, "42+23=", 2) | Out-Null
$d.Content.Find.Execute("30-14=", $true, $false, $false, $false, $false, $true, 1, $false, "2+69=", 2) | Out-Null
$d.Content.Find.Execute("68-3=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("92-63=", $true, $false, $false, $false, $false, $true, 1, $false, "75-74=", 2) | Out-Null
$d.Content.Find.Execute("20-3=", $true, $false, $false, $false, $false, $true, 1, $false, "3+81=", 2) | Out-Null
$d.Content.Find.Execute("46+36=", $true, $false, $false, $false, $false, $true, 1, $false, "41+31=", 2) | Out-Null
$d.Content.Find.Execute("10+79=", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=", 2) | Out-Null
$d.Content.Find.Execute("35+13=", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=", 2) | Out-Null
$d.Content.Find.Execute("74-66=", $true, $false, $false, $false, $false, $true, 1, $false, "19+32=", 2) | Out-Null
$d.Content.Find.Execute("78+0=", $true, $false, $false, $false, $false, $true, 1, $false, "38-31=", 2) | Out-Null
$d.Content.Find.Execute("66-32=", $true, $false, $false, $false, $false, $true, 1, $false, "15+84=", 2) | Out-Null
